# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# for the cryptos worksheet, rows 2-51, per the commit's refreshed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2: 23.368.12 -> 23.358.83; E2: '  +1.44%  ' -> '  +1.34%  '
$ws.Range("D2").Value = '23.358.83'
$ws.Range("E2").Value = '  +1.34%  '

# Row 3: D3: 1.619.96 -> 1.623.43; E3: '  +1.85%  ' -> '  +2.11%  '
$ws.Range("D3").Value = '1.623.43'
$ws.Range("E3").Value = '  +2.11%  '

# Row 4: D4: 0.9959 -> 0.9967; E4: '  -0.64%  ' -> '  -0.61%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9967'
$ws.Range("E4").Value = '  -0.61%  '

# Row 5: D5: 306.04 -> 306.19; E5: '  +1.52%  ' -> '  +1.53%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.19'
$ws.Range("E5").Value = '  +1.53%  '

# Row 6: D6: 0.9968 -> 0.9964; E6: '  -0.51%  ' -> '  -0.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9964'
$ws.Range("E6").Value = '  -0.63%  '

# Row 7: D7: 0.3772 -> 0.3774; E7: '  +0.10%  ' -> '  +0.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3774'
$ws.Range("E7").Value = '  +0.22%  '

# Row 8: D8: 53.28 -> 53.24; E8: '  +5.66%  ' -> '  +4.72%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '53.24'
$ws.Range("E8").Value = '  +4.72%  '

# Row 9: D9: 0.3650 -> 0.3653; E9: '  +1.47%  ' -> '  +1.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3653'
$ws.Range("E9").Value = '  +1.66%  '

# Row 10: D10: 1.277 -> 1.275; E10: '  +4.09%  ' -> '  +3.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.275'
$ws.Range("E10").Value = '  +3.99%  '

# Row 11: D11: 0.08182 -> 0.08176; E11: '  +1.40%  ' -> '  +1.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08176'
$ws.Range("E11").Value = '  +1.45%  '

# Row 12: D12: 0.9958 -> 0.9969; E12: '  -0.68%  ' -> '  -0.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9969'
$ws.Range("E12").Value = '  -0.57%  '

# Row 13: D13: 23.19 -> 23.14; E13: '  +5.27%  ' -> '  +5.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.14'
$ws.Range("E13").Value = '  +5.26%  '

# Row 14: D14: 6.658 -> 6.656; E14: '  +2.71%  ' -> '  +2.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.656'
$ws.Range("E14").Value = '  +2.75%  '

# Row 15: D15: 7.417 -> 7.427; E15: '  +2.13%  ' -> '  +2.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.427'
$ws.Range("E15").Value = '  +2.28%  '

# Row 16: D16: 0.00001253 -> 0.00001255; E16: '  +2.20%  ' -> '  +2.44%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001255'
$ws.Range("E16").Value = '  +2.44%  '

# Row 17: D17: 1.614.55 -> 1.618.61; E17: '  +1.34%  ' -> '  +1.88%  '
$ws.Range("D17").Value = '1.618.61'
$ws.Range("E17").Value = '  +1.88%  '

# Row 18: D18: 94.61 -> 94.71; E18: '  +2.36%  ' -> '  +2.26%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.71'
$ws.Range("E18").Value = '  +2.26%  '

# Row 19: D19: 0.06926 -> 0.06929; E19: '  +1.39%  ' -> '  +1.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06929'
$ws.Range("E19").Value = '  +1.43%  '

# Row 20: D20: 18.34 -> 18.33; E20: '  +2.04%  ' -> '  +2.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.33'
$ws.Range("E20").Value = '  +2.28%  '

# Row 21: D21: 6.575 -> 6.576; E21: '  +1.75%  ' -> '  +1.86%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.576'
$ws.Range("E21").Value = '  +1.86%  '

# Row 22: D22: 1.002 -> 0.9979; E22: '  -0.06%  ' -> '  -0.57%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9979'
$ws.Range("E22").Value = '  -0.57%  '

# Row 23: D23: 12.97 -> 12.96; E23: '  +0.56%  ' -> '  +0.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.96'
$ws.Range("E23").Value = '  +0.78%  '

# Row 24: D24: 23.359.14 -> 23.352.25; E24: '  +1.42%  ' -> '  +1.30%  '
$ws.Range("D24").Value = '23.352.25'
$ws.Range("E24").Value = '  +1.30%  '

# Row 25: D25: 3.129 -> 3.130; E25: '  +11.76%  ' -> '  +11.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.130'
$ws.Range("E25").Value = '  +11.49%  '

# Row 26: D26: 2.420 -> 2.435; E26: '  +1.74%  ' -> '  +2.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.435'
$ws.Range("E26").Value = '  +2.68%  '

# Row 27: E27: '  +2.09%  ' -> '  +2.28%  '
$ws.Range("E27").Value = '  +2.28%  '

# Row 28: D28: 150.59 -> 150.71; E28: '  +1.51%  ' -> '  +1.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '150.71'
$ws.Range("E28").Value = '  +1.57%  '

# Row 29: D29: 5.280 -> 5.279; E29: '  +1.33%  ' -> '  +1.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.279'
$ws.Range("E29").Value = '  +1.06%  '

# Row 30: D30: 136.07 -> 136.15; E30: '  +2.25%  ' -> '  +2.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '136.15'
$ws.Range("E30").Value = '  +2.24%  '

# Row 31: D31: 2.405 -> 2.404; E31: '  +1.70%  ' -> '  +1.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.404'
$ws.Range("E31").Value = '  +1.91%  '

# Row 32: D32: 6.840 -> 6.883; E32: '  +4.16%  ' -> '  +5.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.883'
$ws.Range("E32").Value = '  +5.95%  '

# Row 33: D33: 1.788.93 -> 1.791.60; E33: '  +0.97%  ' -> '  +1.37%  '
$ws.Range("D33").Value = '1.791.60'
$ws.Range("E33").Value = '  +1.37%  '

# Row 34: D34: 0.9665 -> 0.9689; E34: '  +2.03%  ' -> '  +2.79%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9689'
$ws.Range("E34").Value = '  +2.79%  '

# Row 35: E35: '  +3.85%  ' -> '  +4.07%  '
$ws.Range("E35").Value = '  +4.07%  '

# Row 36: D36: 10.39 -> 10.45; E36: '  +2.88%  ' -> '  +3.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.45'
$ws.Range("E36").Value = '  +3.37%  '

# Row 37: D37: 0.07424 -> 0.07436; E37: '  +0.18%  ' -> '  +0.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.07436'
$ws.Range("E37").Value = '  +0.52%  '

# Row 38: D38: 6.198 -> 6.209; E38: '  +2.20%  ' -> '  +2.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.209'
$ws.Range("E38").Value = '  +2.81%  '

# Row 39: D39: 0.2522 -> 0.2528; E39: '  +1.59%  ' -> '  +2.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2528'
$ws.Range("E39").Value = '  +2.00%  '

# Row 40: D40: 0.08802 -> 0.08826; E40: '  +0.11%  ' -> '  +0.64%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08826'
$ws.Range("E40").Value = '  +0.64%  '

# Row 41: D41: 1.410 -> 1.408; E41: '  +3.93%  ' -> '  +4.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.408'
$ws.Range("E41").Value = '  +4.47%  '

# Row 42: D42: 0.7157 -> 0.7152; E42: '  +3.31%  ' -> '  +3.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7152'
$ws.Range("E42").Value = '  +3.50%  '

# Row 43: D43: 12.66 -> 12.67; E43: '  +4.09%  ' -> '  +4.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.67'
$ws.Range("E43").Value = '  +4.42%  '

# Row 44: D44: 16.04 -> 16.02; E44: '  +6.81%  ' -> '  +7.90%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.02'
$ws.Range("E44").Value = '  +7.90%  '

# Row 45: D45: 0.6588 -> 0.6594; E45: '  +1.74%  ' -> '  +2.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6594'
$ws.Range("E45").Value = '  +2.14%  '

# Row 46: D46: 2.347 -> 2.349; E46: '  +4.01%  ' -> '  +4.52%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.349'
$ws.Range("E46").Value = '  +4.52%  '

# Row 47: D47: 0.9956 -> 0.9950; E47: '  -0.51%  ' -> '  -0.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9950'
$ws.Range("E47").Value = '  -0.64%  '

# Row 48: D48: 4.024 -> 4.021; E48: '  +0.62%  ' -> '  +0.59%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.021'
$ws.Range("E48").Value = '  +0.59%  '

# Row 49: D49: 0.08004 -> 0.08012; E49: '  +1.34%  ' -> '  +1.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08012'
$ws.Range("E49").Value = '  +1.51%  '

# Row 50: D50: 131.52 -> 131.50; E50: '  +0.28%  ' -> '  +0.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.50'
$ws.Range("E50").Value = '  +0.18%  '

# Row 51: D51: 1.206 -> 1.210; E51: '  +0.37%  ' -> '  +1.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.210'
$ws.Range("E51").Value = '  +1.09%  '
